$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Java 1.6 (JIT)"
$ws.Range("C1").Value = "Objeck (JIT)"

$ws.Range("A2").Value = 21.063600000000001
$ws.Range("B2").Value = 0.877
$ws.Range("C2").Value = 0.86802199999999996

$ws.Range("A3").Value = 21.0642
$ws.Range("B3").Value = 0.88300000000000001
$ws.Range("C3").Value = 0.86253899999999994

$ws.Range("A4").Value = 21.0655
$ws.Range("B4").Value = 0.93100000000000005
$ws.Range("C4").Value = 0.87425799999999998

$ws.Range("A5").Value = 21.063400000000001
$ws.Range("B5").Value = 0.91200000000000003
$ws.Range("C5").Value = 0.869147

$ws.Rows("14").Delete()

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s = $chart.SeriesCollection().Item(1)
$s.Formula = '=SERIES("Time in Secs",Sheet1!$A$1:$D$1,Sheet1!$A$6:$D$6,1)'
